$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format numeric-looking Price cells as Text so literal strings
# like "1.003" round-trip as text (matching source inline strings)
# instead of being auto-converted to numbers by Excel.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D7:D13").NumberFormat = "@"
$ws.Range("D15:D17").NumberFormat = "@"
$ws.Range("D20:D22").NumberFormat = "@"
$ws.Range("D24:D32").NumberFormat = "@"
$ws.Range("D34:D45").NumberFormat = "@"
$ws.Range("D47:D51").NumberFormat = "@"

# Apply the updated cell values row by row.
$ws.Range("D2").Value = "30.275.15"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "2.095.78"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("D7").Value = "0.5278"
$ws.Range("E7").Value = "  +2.11%  "
$ws.Range("D8").Value = "0.4380"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "55.11"
$ws.Range("E9").Value = "  +2.99%  "
$ws.Range("D10").Value = "0.09347"
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").Value = "1.173"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "24.65"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "8.563"
$ws.Range("E13").Value = "  +5.17%  "
$ws.Range("D14").Value = "2.115.17"
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").Value = "6.859"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("D16").Value = "100.98"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("D17").Value = "0.00001157"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D20").Value = "0.06726"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").Value = "6.381"
$ws.Range("E21").Value = "  +2.83%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("D23").Value = "30.279.50"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("D24").Value = "12.41"
$ws.Range("D25").Value = "2.320"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").Value = "7.028"
$ws.Range("E26").Value = "  +10.08%  "
$ws.Range("D27").Value = "21.79"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").Value = "162.55"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").Value = "2.509"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").Value = "133.71"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "1.133"
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").Value = "1.673"
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").Value = "6.239"
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("D35").Value = "3.911"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("D36").Value = "10.07"
$ws.Range("E36").Value = "  -3.37%  "
$ws.Range("D37").Value = "0.02616"
$ws.Range("E37").Value = "  +1.56%  "
$ws.Range("D38").Value = "0.06756"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").Value = "12.61"
$ws.Range("E39").Value = "  +1.47%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.6961"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.338"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "0.2211"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "0.6774"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").Value = "14.35"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").Value = "2.340"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").Value = "1.303"
$ws.Range("E47").Value = "  +8.53%  "
$ws.Range("D48").Value = "3.639"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("D49").Value = "0.00000000351"
$ws.Range("E49").Value = "  -2.21%  "
$ws.Range("D50").Value = "1.208"
$ws.Range("E50").Value = "  +5.67%  "
$ws.Range("D51").Value = "1.212"
$ws.Range("E51").Value = "  -0.27%  "
